$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness column (C2:C252) from 7569 to 7293
$ws.Range("C2:C252").Value = 7293
